$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.463.03"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.574.44"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.26"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("E7").Value = "  -0.73%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.94"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3403"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07556"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.140"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.995"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.950"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.566.43"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("E17").Value = "  -1.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.10"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06757"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.273"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.463.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.338"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.596"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.60"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.005"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.86"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.742.82"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.054"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.124"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.983"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.819"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08412"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.387"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02467"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06541"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.461"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.29"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6262"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.10%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.95"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.812"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5830"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.088"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.12"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.225"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07325"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.08%  "
